# GSCHED-578: Fixing spreadsheet columns and test cases for guiders.
#
# 1. Rename instrument/guider labels:
#    "F2"       -> "Flamingos2"
#    "GMOS OI"  -> "GMOS OIWFS"
#    "F2 OI"    -> "FII OIWFS"
# 2. Fix the style applied to column J (Port 5) for rows 424-547 on the
#    "GS" sheet so it matches the style already used elsewhere (style
#    index 6, which differs from style index 5 only by readingOrder).

$wb = $excel.ActiveWorkbook

# --- Part 1: rename shared strings wherever they occur on any sheet ---
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        switch ($cell.Value) {
            "F2"      { $cell.Value = "Flamingos2" }
            "GMOS OI" { $cell.Value = "GMOS OIWFS" }
            "F2 OI"   { $cell.Value = "FII OIWFS" }
        }
    }
}

# --- Part 2: fix the style of J424:J547 on the GS sheet ---
$gs = $wb.Worksheets.Item("GS")

# Find a cell that already carries the "correct" (style index 6) look so
# we can copy its style onto the target range instead of hard-coding a
# style name.
$template = $gs.Range("K493")

$target = $gs.Range("J424:J547")
$target.Style = $template.Style
